# Add a "status_label" column (string version of the "statut" emoji column)
# right after column A, shifting the existing columns (NCTId..results) one
# place to the right, and populate the new column with the French label
# "rouge" for every data row. Also fix a data mix-up where rows 3 and 4 had
# their NCTId / clinical_trial_title / acronym swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B; this shifts B:I -> C:J and
# carries the header style along for row 1.
$ws.Columns.Item(2).Insert()

# New header + values for the inserted "status_label" column.
$ws.Cells.Item(1, 2).Value = "status_label"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "rouge"
}

# Correct the NCTId / clinical_trial_title / acronym values that were
# swapped between rows 3 and 4 (column letters after the insert: C, F, G).
$ws.Cells.Item(3, 3).Value = "NCT02822209"
$ws.Cells.Item(3, 6).Value = "Evaluation of the Impact of a Coordinating Nurse in a Personalized Care Program on Quality of Care, Coordination of the Actors and on Quality of Life for Patients With Lung Cancer. A French Randomized Monocentric Prospective Study"
$ws.Cells.Item(3, 7).Value = "EVIDEC"

$ws.Cells.Item(4, 3).Value = "NCT02821637"
$ws.Cells.Item(4, 6).Value = "Impact of an Effort Rehabilitation Program for Overweight or Obese Children and Teens on Quality of Life and Wellbeing: A French Monocentric Prospective Study"
$ws.Cells.Item(4, 7).Value = "PRESEVAL"
